$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 62, shifting rows 62..133 down to 63..134.
$ws.Rows.Item(62).Insert()

# Populate the newly-inserted row 62 with the new record's data.
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 44740
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = 100112032
$ws.Cells.Item(62, 7).Value = "Zapallo italiano"
$ws.Cells.Item(62, 8).Value = "Sin especificar"
$ws.Cells.Item(62, 9).Value = "Primera"
$ws.Cells.Item(62, 10).Value = 250
$ws.Cells.Item(62, 11).Value = 14000
$ws.Cells.Item(62, 12).Value = 15000
$ws.Cells.Item(62, 13).Value = 14600
$ws.Cells.Item(62, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(62, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(62, 16).Value = 243
$ws.Cells.Item(62, 17).Value = 60
$ws.Cells.Item(62, 18).Value = "Hortaliza"
